$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Pearl Morisette genre changed from "Canadian" to "Contemporary"
$ws.Range("D12").Value = "Contemporary"

# --- Row 26: Miss Fu in Shengdu - "Don't Miss" and Genre updated
$ws.Range("C26").Value = "Chongqing Chicken + Spicy Potatoes; Sichuan goodness"
$ws.Range("D26").Value = "Chinese"

# --- New rows 33-39 appended
$ws.Range("A33").Value = "Yummy Chinese"
$ws.Range("B33").Value = "Leslieville"
$ws.Range("C33").Value = "Popcorn chicken is wildly addictive, very reasonable price, best takeout Cantonese style Chinese"
$ws.Range("D33").Value = "Chinese"
$ws.Range("E33").Value = 43.670188149640403
$ws.Range("F33").Value = -79.336298702573899

$ws.Range("A34").Value = "Gio Rana's Really Really Nice"
$ws.Range("B34").Value = "Leslieville"
$ws.Range("C34").Value = "Giant Meatball, really anything? Killer Italian. "
$ws.Range("D34").Value = "Italian"
$ws.Range("E34").Value = 43.663396670633396
$ws.Range("F34").Value = -79.330474179261202

$ws.Range("A35").Value = "Enoteca Sociale"
$ws.Range("B35").Value = "Little Portugal"
$ws.Range("C35").Value = "Fresh pastas, mocktails were great, stunning meal"
$ws.Range("D35").Value = "Italian"
$ws.Range("E35").Value = 43.649780159325701
$ws.Range("F35").Value = -79.425617316722807

$ws.Range("A36").Value = "MIMI Chinese"
$ws.Range("B36").Value = "Yorkville"
$ws.Range("C36").Value = "Shrimp Toast was stunning, Hidden Crispy Chicken, very playful and flavorful and everything good in the world"
$ws.Range("D36").Value = "Chinese"
$ws.Range("E36").Value = 43.674636892456697
$ws.Range("F36").Value = -79.3981557604356

$ws.Range("A37").Value = "Que Ling Vietnamese Cuisine"
$ws.Range("B37").Value = "Leslieville"
$ws.Range("C37").Value = "Pho + noodle soup, cash only"
$ws.Range("D37").Value = "Vietnamese"
$ws.Range("E37").Value = 43.665957195598601
$ws.Range("F37").Value = -79.349426857565504

$ws.Range("A38").Value = "Ocha's "
$ws.Range("B38").Value = "Leslieville"
$ws.Range("C38").Value = "Jerk chicken, festival; have heard great things on their oxtail but have yet to try"
$ws.Range("D38").Value = "Caribbean"
$ws.Range("E38").Value = 43.661691129281799
$ws.Range("F38").Value = -79.338724543548807

$ws.Range("A39").Value = "Taqueria el pastorcito"
$ws.Range("B39").Value = "Allan Gardens ish?"
$ws.Range("C39").Value = "al pastor! "
$ws.Range("D39").Value = "Mexican"
$ws.Range("E39").Value = 43.664608478470498
$ws.Range("F39").Value = -79.384586629285295

# --- Update sheet view to match (selection)
$ws.Range("E40").Select()
